$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 79244

# Row 3
$ws.Range("A3").Value = 131067787
$ws.Range("B3").Value = 57884
$ws.Range("E3").Value = 100109
$ws.Range("F3").Value = "Tretåig hackspett"
$ws.Range("G3").Value = "Picoides tridactylus"
$ws.Range("H3").Value = "(Linnaeus, 1758)"
$ws.Range("J3").Value = ""
$ws.Range("L3").Value = ""
$ws.Range("M3").Value = "äldre spår"
$ws.Range("Q3").Value = 466335
$ws.Range("R3").Value = 7046445
$ws.Range("AC3").Value = "Ringhack, äldre, på gran."
$ws.Range("AF3").Value = ""
$ws.Range("AM3").Value = ""
$ws.Range("AO3").Value = "Picea abies"

# Row 4
$ws.Range("A4").Value = 131067798
$ws.Range("Q4").Value = 466279
$ws.Range("R4").Value = 7046403
$ws.Range("AH4").Value = "Granskog"

# Row 5
$ws.Range("A5").Value = 131067788
$ws.Range("Q5").Value = 466325
$ws.Range("R5").Value = 7046442
$ws.Range("AH5").Value = "Barrskog"

# Row 6
$ws.Range("A6").Value = 131067035
$ws.Range("B6").Value = 79244
$ws.Range("E6").Value = 6425
$ws.Range("F6").Value = "Garnlav"
$ws.Range("G6").Value = "Alectoria sarmentosa"
$ws.Range("H6").Value = "(Ach.) Ach."
$ws.Range("J6").Value = ""
$ws.Range("L6").Value = ""
$ws.Range("M6").Value = ""
$ws.Range("Q6").Value = 466172
$ws.Range("R6").Value = 7046340
$ws.Range("AC6").Value = "Långväxta bålar på gran."
$ws.Range("AF6").Value = ""
$ws.Range("AM6").Value = "Gren på levande träd"
$ws.Range("AO6").Value = "Branch on living tree # Picea abies"

# Row 10
$ws.Range("A10").Value = 131067786
$ws.Range("B10").Value = 57884
$ws.Range("E10").Value = 100109
$ws.Range("F10").Value = "Tretåig hackspett"
$ws.Range("G10").Value = "Picoides tridactylus"
$ws.Range("H10").Value = "(Linnaeus, 1758)"
$ws.Range("J10").Value = ""
$ws.Range("L10").Value = ""
$ws.Range("M10").Value = "äldre spår"
$ws.Range("Q10").Value = 466366
$ws.Range("R10").Value = 7046466
$ws.Range("AC10").Value = "Ringhack, äldre, på gran."
$ws.Range("AF10").Value = ""
$ws.Range("AM10").Value = ""
$ws.Range("AO10").Value = "Picea abies"

# Row 11
$ws.Range("A11").Value = 131067792
$ws.Range("M11").Value = "färska spår"
$ws.Range("Q11").Value = 466356
$ws.Range("R11").Value = 7046460
$ws.Range("AC11").Value = "Ringhack, färska, på gran."
$ws.Range("AM11").Value = "Trädstam på levande träd"
$ws.Range("AO11").Value = "Stem on living tree # Picea abies"

# Row 12
$ws.Range("A12").Value = 131067781
$ws.Range("M12").Value = "äldre spår"
$ws.Range("Q12").Value = 466204
$ws.Range("R12").Value = 7046448
$ws.Range("AC12").Value = "Ringhack, äldre, på gran."
$ws.Range("AH12").Value = "Granskog"
$ws.Range("AM12").Value = ""
$ws.Range("AO12").Value = "Picea abies"

# Row 13
$ws.Range("A13").Value = 131067030
$ws.Range("B13").Value = 79244
$ws.Range("E13").Value = 6425
$ws.Range("F13").Value = "Garnlav"
$ws.Range("G13").Value = "Alectoria sarmentosa"
$ws.Range("H13").Value = "(Ach.) Ach."
$ws.Range("J13").Value = ""
$ws.Range("L13").Value = ""
$ws.Range("M13").Value = ""
$ws.Range("Q13").Value = 466302
$ws.Range("R13").Value = 7046517
$ws.Range("AC13").Value = ""
$ws.Range("AF13").Value = ""
$ws.Range("AH13").Value = "Barrskog"
$ws.Range("AM13").Value = "Gren på levande träd"
$ws.Range("AO13").Value = "Branch on living tree # Picea abies"

# Row 26
$ws.Range("B26").Value = 79244

# Row 27
$ws.Range("A27").Value = 131067038
$ws.Range("B27").Value = 79244
$ws.Range("E27").Value = 6425
$ws.Range("F27").Value = "Garnlav"
$ws.Range("G27").Value = "Alectoria sarmentosa"
$ws.Range("H27").Value = "(Ach.) Ach."
$ws.Range("J27").Value = ""
$ws.Range("L27").Value = ""
$ws.Range("M27").Value = ""
$ws.Range("Q27").Value = 466024
$ws.Range("R27").Value = 7046276
$ws.Range("AC27").Value = "Enstaka bålar på gran."
$ws.Range("AF27").Value = ""
$ws.Range("AM27").Value = "Gren på levande träd"
$ws.Range("AO27").Value = "Branch on living tree # Picea abies"

# Row 28
$ws.Range("A28").Value = 131067797
$ws.Range("B28").Value = 57884
$ws.Range("E28").Value = 100109
$ws.Range("F28").Value = "Tretåig hackspett"
$ws.Range("G28").Value = "Picoides tridactylus"
$ws.Range("H28").Value = "(Linnaeus, 1758)"
$ws.Range("J28").Value = ""
$ws.Range("L28").Value = ""
$ws.Range("M28").Value = "äldre spår"
$ws.Range("Q28").Value = 466283
$ws.Range("R28").Value = 7046407
$ws.Range("AC28").Value = "Ringhack, äldre, på gran."
$ws.Range("AF28").Value = ""
$ws.Range("AM28").Value = ""
$ws.Range("AO28").Value = "Picea abies"

# Row 30
$ws.Range("B30").Value = 79244

# Row 31
$ws.Range("B31").Value = 79244

# Row 35
$ws.Range("B35").Value = 79244

# Row 36
$ws.Range("B36").Value = 79244

# Row 39
$ws.Range("B39").Value = 79244

# Row 40
$ws.Range("B40").Value = 79244
